# Update the build timestamp embedded in the version string across the
# workbook: "February 03 2026 17.29.55 EST" -> "February 03 2026 18.05.36 EST"

$wb = $excel.ActiveWorkbook

$oldStamp = "February 03 2026 17.29.55 EST"
$newStamp = "February 03 2026 18.05.36 EST"

$aboutSheet = $wb.Worksheets.Item("About")
$dataSheet = $wb.Worksheets.Item("Boundaries and methane sources")

# --- "About" sheet: A2 (version banner) and A6 (recommended citation) ---
$aboutSheet.Range("A2").Value = $aboutSheet.Range("A2").Value().Replace($oldStamp, $newStamp)
$aboutSheet.Range("A6").Value = $aboutSheet.Range("A6").Value().Replace($oldStamp, $newStamp)

# --- "Boundaries and methane sources" sheet: column S, rows 2-66 ---
for ($row = 2; $row -le 66; $row++) {
    $cell = $dataSheet.Cells.Item($row, 19)  # column S = 19
    $val = $cell.Value()
    if ($val -ne $null -and $val.ToString().Contains($oldStamp)) {
        $cell.Value = $val.ToString().Replace($oldStamp, $newStamp)
    }
}
